$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.866.20'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.628.77'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').Value = '214.73'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').Value = '28.64'
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = '0.258'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '0.0609'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').Value = '0.0900'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('D12').Value = '1.860.66'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '1.630.09'
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').Value = '0.566'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = '9.36'
$ws.Range('E15').Value = '  +5.40%  '
$ws.Range('D16').Value = '29.879.81'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '3.84'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '65.00'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '241.14'
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('D20').Value = '0.0₃0703'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').Value = '9.79'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('E24').Value = '  +3.00%  '
$ws.Range('D25').Value = '157.60'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').Value = '15.48'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').Value = '6.56'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').Value = '0.0489'
$ws.Range('E30').Value = '  +1.39%  '
$ws.Range('D31').Value = '1.11'
$ws.Range('E31').Value = '  +3.35%  '
$ws.Range('D32').Value = '3.37'
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('D33').Value = '3.17'
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('D34').Value = '1.427.25'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').Value = '1.68'
$ws.Range('E35').Value = '  +2.95%  '
$ws.Range('D36').Value = '1.02'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('E37').Value = '  -4.62%  '
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').Value = '0.0171'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').Value = '75.01'
$ws.Range('E40').Value = '  +7.92%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '1.98'
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('D43').Value = '0.830'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '0.0498'
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '5.33'
$ws.Range('E47').Value = '  -2.09%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.767.87'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').Value = '48.92'
$ws.Range('E49').Value = '  -8.87%  '
$ws.Range('D50').Value = '91.50'
$ws.Range('E50').Value = '  +4.19%  '
$ws.Range('E51').Value = '  +8.57%  '
